$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (Bmw m8): mark as sold, add sale price, highlight car name in red ---
$ws.Range("F5").Value = 15700000
$ws.Range("B5").Font.Bold = $true
$ws.Range("B5").Font.Italic = $true
$ws.Range("B5").Font.Color = 255

# --- Row 15 (MB S63): mark as sold, add sale price, highlight car name in red ---
$ws.Range("F15").Value = 17300000
$ws.Range("B15").Font.Bold = $true
$ws.Range("B15").Font.Italic = $true
$ws.Range("B15").Font.Color = 255

# --- Row 16 (MB G65): shorten note, mark as sold, add sale price, highlight car name, add comment ---
$ws.Range("D16").Value = "номер + 4 кк"
$ws.Range("F16").Value = 22000000
$ws.Range("B16").Font.Bold = $true
$ws.Range("B16").Font.Italic = $true
$ws.Range("B16").Font.Color = 255
$ws.Range("G16").Value = "Кирилл угадал"
$ws.Range("G16").Font.Bold = $true

# --- Highlight the "Total cars" sum cell to match the "Total sum" styling ---
$ws.Range("C23").Interior.Color = 49407

# --- Move the active selection ---
$ws.Range("L11").Select()
